$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column letters (in order) corresponding to the values in each row of $data below.
$cols = @("B", "C", "D", "E", "G", "H", "L", "O")

$data = @(
    @(17.15386298278785, 11.52759821886406, 5.978516773887472, 11.86158874863263, 27.47276373009179, 13.91327082812218, 9.694535064815888, 21.01690839712343),
    @(16.47104559941554, 11.31770863479185, 5.85862166476447, 11.92341555060806, 27.54629610950944, 13.97730051633524, 9.659935832854574, 21.11579530358568),
    @(16.03784181814913, 11.18696006769778, 5.785545473284466, 11.96383183875734, 27.60500295813436, 14.01980993267216, 9.640236141087627, 21.18318059383411),
    @(15.85803954160074, 11.13326258255971, 5.755948342617581, 11.98091908160068, 27.63231001128509, 14.03793469294474, 9.632602210609001, 21.21230926571885),
    @(15.82799369640914, 11.12432265632084, 5.751046114915264, 11.98379369083317, 27.63704788385603, 14.04099267849174, 9.631358543071372, 21.21724659747058),
    @(16.03542983648963, 11.18623749627664, 5.785145517614804, 11.96405978380359, 27.605357566657, 14.02005112503604, 9.640131585554904, 21.1835666880906),
    @(16.92145826022709, 11.4556501153161, 5.937092460967837, 11.88239728175875, 27.49528917589865, 13.9346842177091, 9.682288136825989, 21.04961571903221),
    @(18.53912894121093, 11.96675920005878, 6.237347804247364, 11.74172710884976, 27.38803344198093, 13.79270015506236, 9.776952041325732, 20.84020683815409),
    @(19.6441575762748, 12.32882052974822, 6.456749503359978, 11.65023961286593, 27.37661233936532, 13.70397689195829, 9.853443762280611, 20.71930609393025),
    @(20.12699796059694, 12.49003266606925, 6.55578511980653, 11.61119475989391, 27.38623321135563, 13.667020856493, 9.889663005737145, 20.67155835867053),
    @(20.30686521465343, 12.55053264267706, 6.593136159200882, 11.5967796087685, 27.39201692521089, 13.6535179254794, 9.903575072751163, 20.65452785511143),
    @(20.26826158379099, 12.53752793844371, 6.585099310341944, 11.59986769644003, 27.39067597995895, 13.65640413201296, 9.900570245165564, 20.65814883363608),
    @(20.14185597811175, 12.4950212603619, 6.558861275449579, 11.61000139404702, 27.38666609813653, 13.66590010038172, 9.890803665861462, 20.67013615056201),
    @(20.06403832387124, 12.46891208049567, 6.542768810392055, 11.61625680691027, 27.38448891903468, 13.671780722839, 9.884846719050042, 20.67761576740787),
    @(19.6121928232802, 12.31821072960399, 6.450258242365285, 11.65284308754111, 27.37628266591842, 13.70646073098797, 9.851104705055139, 20.72257310337008),
    @(19.32982655446643, 12.22483286150403, 6.393278789908233, 11.67594683163445, 27.37505141674434, 13.7286091750995, 9.830763905092974, 20.75201573986475),
    @(19.16555525971405, 12.17079830537372, 6.360435410511496, 11.68947768495221, 27.37573745348623, 13.74166871727746, 9.819198829829148, 20.76963237330268),
    @(19.10962004266649, 12.1524485605891, 6.349304408269661, 11.69410058870372, 27.37620883153249, 13.74614541413998, 9.815306417854606, 20.77571395289792),
    @(19.36007853271127, 12.23480718095416, 6.399351928191047, 11.67346233105658, 27.37503810565901, 13.72621826708098, 9.832915363324572, 20.74881087755965),
    @(20.17906595417206, 12.50752171175024, 6.566572449535083, 11.60701482972606, 27.38778574384932, 13.66309754918567, 9.893667075521904, 20.6665866124984),
    @(20.69694329966601, 12.68254565729927, 6.674959862933151, 11.56574600667323, 27.40859557616617, 13.62471035238773, 9.934514021607743, 20.61897503303904),
    @(20.4221679297748, 12.58944019067717, 6.617206732255247, 11.58757436198741, 27.3963448460913, 13.64493545638196, 9.912611432954966, 20.64382310030879),
    @(19.34640763163852, 12.23029887900675, 6.396606524446087, 11.6745848006464, 27.37503978245831, 13.72729818090835, 9.831942286262892, 20.75025764808599),
    @(18.11550396807645, 11.83064774098864, 6.15614918644165, 11.77769853098831, 27.4052896729237, 13.82837936172389, 9.750095560380208, 20.89110452827137)
)

# Data rows start at worksheet row 2 (row 1 holds headers).
$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowValues[$j]
    }
}
